# Add a new "Sheet3" worksheet after the existing sheets and populate it
# with a small admin-user search test fixture, then update the view
# selections so Sheet3 becomes the active tab (matching "added test for
# search admin user").

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing worksheet so it lands at
# the end of the tab order (Sheet1, Sheet2, Sheet3).
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Sheet3"

# Populate the new test data.
$ws3.Range("A1").Value = "UserRole"
$ws3.Range("B1").Value = "Staus"
$ws3.Range("A2").Value = "Admin"
$ws3.Range("B2").Value = "Enabled"

# Sheet2 keeps a selection but is no longer the active/selected tab.
$ws2 = $wb.Worksheets.Item("Sheet2")
[void]$ws2.Select()
[void]$ws2.Range("G2").Select()

# Sheet3 becomes the active tab with its own selection.
[void]$ws3.Select()
[void]$ws3.Range("K17").Select()
